# Updated cryptos list on Sat May 18 07:10:03 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in columns D/E (and the B/C text in the swapped rows) are
# stored as text in the workbook, even when they look numeric (e.g. "579.48"),
# or percentages with padding (e.g. "  +1.17%  "). Force text format first so
# Excel does not silently convert these into numeric cells on write.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.885.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.097.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.36%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.77%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.094.47"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.41"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.28%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.42"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.87%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.609.68"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.859.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.100.26"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.04"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.04"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.24"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.30%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.988"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.98"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.93%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.14"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.317"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.41%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.845.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.72%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "383.71"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.72"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.97"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.35%  "
